$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (Hydrogen): D3 loses its numeric value and becomes blank
$ws.Range("D3").ClearContents()

# Row 4 (Methanol): C4 value changes from 17694.8045882471 to 0
$ws.Range("C4").Value = 0

# Row 5 (Ammonia): C5 value changes from 72228.60202626928 to 1255.355988682103
$ws.Range("C5").Value = 1255.355988682103

# Row 7: rename "Other" -> "Biogas", update D7 value
$ws.Range("A7").Value = "Biogas"
$ws.Range("D7").Value = 235.728752017606

# Row 8 (new row): copy formatting (and blank B/C cells) from row 7, then set the new "Other" entry
$ws.Range("A7:D7").Copy($ws.Range("A8:D8"))
$ws.Range("A8").Value = "Other"
$ws.Range("D8").Value = 187.8568869578981
